# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update the DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Row = 8;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 24; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 36; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 38; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 39; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 40; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 41; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 48; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 50; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 59; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 67; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 80; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 88; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
